$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Absent" column (H) values for the consolidated report.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
